# Updated cryptos list values (prices & volume deltas) per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.572.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.139.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "640.81"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +10.96%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.23%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.136.52"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.720"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.63"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.441.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.720.50"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.154.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.46"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.56%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.90%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.05"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.11"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.33%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.44"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +15.99%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +33.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.89"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "518.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.14"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.31"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("B40").Value = "Binance-PegBSC-USD"
$ws.Range("C40").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -10.48%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.423"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.73%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.20"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0860"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +41.96%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.707"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +13.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "150.63"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.74"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.58"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +9.38%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.31%  "
